$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlCenter
$xlCenter = -4108

# --- Header row (row 1): repeat the B1/C1/D1 headers in I1/J1/K1 ---
$ws.Range("I1").Value = "Serializalble"
$ws.Range("J1").Value = "Remote"
$ws.Range("K1").Value = "UnicastRemoteObject"

# --- New class/interface names in column H (rows 2-8) ---
# NOTE: the order these string values are first written controls the order
# they land in the shared-string table, so they are written in the same
# sequence the original author used (IRekening, IKlant, Money, Balie,
# Bankiersessie, IBalie, IBankiersessie) rather than strict row order.
$ws.Range("H3").Value = "IRekening"
$ws.Range("H2").Value = "IKlant"
$ws.Range("H4").Value = "Money"
$ws.Range("H5").Value = "Balie"
$ws.Range("H7").Value = "Bankiersessie"
$ws.Range("H6").Value = "IBalie"
$ws.Range("H8").Value = "IBankiersessie"

# --- Marks ("x") in columns I/J/K for rows 2-8 ---
$ws.Range("I2").Value = "x"
$ws.Range("I3").Value = "x"
$ws.Range("I4").Value = "x"
$ws.Range("K5").Value = "x"
$ws.Range("J6").Value = "x"
$ws.Range("K7").Value = "x"
$ws.Range("J8").Value = "x"

# --- Apply the centered style (style index 1, same style already used by the
# column-D data cells) to every I/J/K cell that needs it in rows 2-8. Cells
# that already got a value above pick this style up too; the rest become
# styled-but-empty cells. Row 6 intentionally gets no K6 cell (matches the
# source workbook, where that cell was never touched).
$ws.Range("I2:K4").HorizontalAlignment = $xlCenter
$ws.Range("I5:K5").HorizontalAlignment = $xlCenter
$ws.Range("I6:J6").HorizontalAlignment = $xlCenter
$ws.Range("I7:K7").HorizontalAlignment = $xlCenter
$ws.Range("I8:K8").HorizontalAlignment = $xlCenter

# --- Rows 9-33: blank styled cells in I, J, K only (no H cell there) ---
$ws.Range("I9:K33").HorizontalAlignment = $xlCenter

# --- Column widths for H, I, K ---
# (the engine snaps ColumnWidth to a 1/6-character pixel grid, so these are
# the inputs that land closest to the target stored widths of
# 14 / 11.5703125 / 20.28515625)
$ws.Range("H1").ColumnWidth = 13.166666666666668
$ws.Range("I1").ColumnWidth = 10.6665
$ws.Range("K1").ColumnWidth = 19.41825

# --- View state: selection on K17 (scrolled-to column not supported by this host) ---
$ws.Range("K17").Select()
